$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40780
$ws.Range("J3").Value = 40780
$ws.Range("L3").Value = 40780
$ws.Range("N3").Value = -41008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 648.6667
$ws.Range("I32").Value = 495
$ws.Range("J32").Value = 771.6
$ws.Range("K32").Value = 495
$ws.Range("L32").Value = 771.6
$ws.Range("M32").Value = -169
$ws.Range("N32").Value = -1423.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2774.875
$ws.Range("J88").Value = 2766.3333
$ws.Range("L88").Value = 2766.3333
$ws.Range("N88").Value = -3578.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2774.875
$ws.Range("J91").Value = 2766.3333
$ws.Range("L91").Value = 2766.3333
$ws.Range("N91").Value = -5574.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 40780
$ws.Range("J102").Value = 40780
$ws.Range("L102").Value = 40780
$ws.Range("N102").Value = -47270

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5681.6875
$ws.Range("I116").Value = 6391.5454
$ws.Range("J116").Value = 4120
$ws.Range("K116").Value = 6391.5454
$ws.Range("L116").Value = 4120
$ws.Range("M116").Value = -2949.5454
$ws.Range("N116").Value = -11004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1872.5
$ws.Range("I131").Value = 495
$ws.Range("J131").Value = 3250
$ws.Range("K131").Value = 1485
$ws.Range("L131").Value = 9750
$ws.Range("M131").Value = 3555
$ws.Range("N131").Value = -19830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 865.0833
$ws.Range("I2").Value = 722.625
$ws.Range("J2").Value = 2004.75
$ws.Range("K2").Value = 722.625
$ws.Range("L2").Value = 2004.75
$ws.Range("M2").Value = -609.625
$ws.Range("N2").Value = -2230.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3584.0605
$ws.Range("I45").Value = 4473.48
$ws.Range("J45").Value = 804.625
$ws.Range("K45").Value = 4473.48
$ws.Range("L45").Value = 804.625
$ws.Range("M45").Value = -4096.48
$ws.Range("N45").Value = -1558.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 364082.3
$ws.Range("I61").Value = 10120.929
$ws.Range("K61").Value = 10120.929
$ws.Range("M61").Value = -9908.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 865.0833
$ws.Range("I116").Value = 722.625
$ws.Range("J116").Value = 2004.75
$ws.Range("K116").Value = 722.625
$ws.Range("L116").Value = 2004.75
$ws.Range("M116").Value = 1571.375
$ws.Range("N116").Value = -6592.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 364082.3
$ws.Range("I136").Value = 10120.929
$ws.Range("K136").Value = 30362.787
$ws.Range("M136").Value = -27812.787

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 865.0833
$ws.Range("I3").Value = 722.625
$ws.Range("J3").Value = 2004.75
$ws.Range("K3").Value = 722.625
$ws.Range("L3").Value = 2004.75
$ws.Range("M3").Value = -608.625
$ws.Range("N3").Value = -2232.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 754.0909
$ws.Range("I11").Value = 500
$ws.Range("K11").Value = 500
$ws.Range("M11").Value = -360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1198.56
$ws.Range("I107").Value = 1180.5
$ws.Range("K107").Value = 1180.5
$ws.Range("M107").Value = 739.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 386304.34
$ws.Range("I58").Value = 1759.3334
$ws.Range("J58").Value = 910683.9399999999
$ws.Range("K58").Value = 1759.3334
$ws.Range("L58").Value = 910683.9399999999
$ws.Range("M58").Value = -1556.3334
$ws.Range("N58").Value = -911089.9399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13096.889
$ws.Range("I99").Value = 14509
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 14509
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -13011
$ws.Range("N99").Value = -4796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3788.75
$ws.Range("I105").Value = 3841.3333
$ws.Range("K105").Value = 3841.3333
$ws.Range("M105").Value = -2094.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 760.03705
$ws.Range("I107").Value = 577
$ws.Range("J107").Value = 1026.2727
$ws.Range("K107").Value = 577
$ws.Range("L107").Value = 1026.2727
$ws.Range("M107").Value = 1343
$ws.Range("N107").Value = -4866.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13096.889
$ws.Range("I126").Value = 14509
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 43527
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -41057
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 246814.56
$ws.Range("I134").Value = 2969.4324
$ws.Range("J134").Value = 2502382
$ws.Range("K134").Value = 8908.297200000001
$ws.Range("L134").Value = 7507146
$ws.Range("M134").Value = -6373.297200000001
$ws.Range("N134").Value = -7512216

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 386304.34
$ws.Range("I136").Value = 1759.3334
$ws.Range("J136").Value = 910683.9399999999
$ws.Range("K136").Value = 5278.0002
$ws.Range("L136").Value = 2732051.82
$ws.Range("M136").Value = -2728.0002
$ws.Range("N136").Value = -2737151.82

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 755590.6
$ws.Range("I132").Value = 540.8
$ws.Range("J132").Value = 935364.4399999999
$ws.Range("K132").Value = 4867.2
$ws.Range("L132").Value = 8418279.959999999
$ws.Range("M132").Value = -2337.2
$ws.Range("N132").Value = -8423339.959999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1582.1724
$ws.Range("I102").Value = 1477.1482
$ws.Range("K102").Value = 1477.1482
$ws.Range("M102").Value = 144.8517999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1454.1154
$ws.Range("I107").Value = 595.8570999999999
$ws.Range("J107").Value = 2455.4167
$ws.Range("K107").Value = 595.8570999999999
$ws.Range("L107").Value = 2455.4167
$ws.Range("M107").Value = 1324.1429
$ws.Range("N107").Value = -6295.4167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 38031364
$ws.Range("I122").Value = 50707484
$ws.Range("J122").Value = 2995.7144
$ws.Range("K122").Value = 152122452
$ws.Range("L122").Value = 8987.143199999999
$ws.Range("M122").Value = -152120002
$ws.Range("N122").Value = -13887.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5231.5776
$ws.Range("I132").Value = 6048.2334
$ws.Range("J132").Value = 3598.2666
$ws.Range("K132").Value = 18144.7002
$ws.Range("L132").Value = 10794.7998
$ws.Range("M132").Value = -15614.7002
$ws.Range("N132").Value = -15854.7998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 45457376
$ws.Range("I40").Value = 58826380
$ws.Range("J40").Value = 2759.8
$ws.Range("K40").Value = 58826380
$ws.Range("L40").Value = 2759.8
$ws.Range("M40").Value = -58826244
$ws.Range("N40").Value = -3031.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 948.8333
$ws.Range("I46").Value = 600.5
$ws.Range("J46").Value = 1123
$ws.Range("K46").Value = 600.5
$ws.Range("L46").Value = 1123
$ws.Range("M46").Value = -412.5
$ws.Range("N46").Value = -1499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1233.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3578978.5
$ws.Range("I122").Value = 4770304.5
$ws.Range("K122").Value = 14310913.5
$ws.Range("M122").Value = -14308463.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9265335
$ws.Range("I132").Value = 13340091
$ws.Range("J132").Value = 4524
$ws.Range("K132").Value = 40020273
$ws.Range("L132").Value = 13572
$ws.Range("M132").Value = -40017743
$ws.Range("N132").Value = -18632

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1229.9016
$ws.Range("I132").Value = 811.617
$ws.Range("J132").Value = 2634.1428
$ws.Range("K132").Value = 2434.851
$ws.Range("L132").Value = 7902.428400000001
$ws.Range("M132").Value = 95.14900000000034
$ws.Range("N132").Value = -12962.4284
